$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update Test Case ID
$ws.Range("B1").Value = "TC-41"

# 2. Update test data: Password field text, then clear the remaining test-data rows
$ws.Range("E11").Value = "Contraseña"
$ws.Range("E12").Value = ""
$ws.Range("E13").Value = ""
$ws.Range("E14").Value = ""
$ws.Range("E15").Value = ""
$ws.Range("E16").Value = ""
$ws.Range("E17").Value = ""
$ws.Range("E18").Value = ""

# 3. Clear step 4 details (text moved back into the form, paso 2 note removed)
$ws.Range("B29").Value = ""
